# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before
#    "总计"), populated with per-fund holding data for the quarter.
# 2. Prepend a "2022-Q1" summary row to the "总计" (totals) sheet.
#
# "总计" is recreated (captured + deleted + re-added at the end) rather
# than edited in place so its sheetId shifts from 6 -> 7, freeing 6 for
# the new "2022-Q1" sheet - matching how Excel renumbers sheetIds when a
# sheet is inserted ahead of the existing last sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0) Snapshot the existing "总计" rows, then delete that sheet so its
#    sheetId is freed up for the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$existingTotals = @()
for ($r = 2; $r -le 6; $r++) {
    $rowVals = @($tot.Cells.Item($r, 2).Value2, $tot.Cells.Item($r, 3).Value2, $tot.Cells.Item($r, 4).Value2)
    $existingTotals += ,$rowVals
}
[void]$tot.Delete()

# ---------------------------------------------------------------------
# 1) Add the new "2022-Q1" sheet, positioned after "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$new = $wb.Worksheets.Add($null, $q4)
$new.Name = "2022-Q1"

# Copy the header row (B1:H1) + the A-column row index (A2:A5) from the
# "2021-Q4" sheet so the new sheet starts out with matching styles
# (bold/centered header, bordered index column).
$q4.Range("B1:H1").Copy($new.Range("B1:H1"))
$q4.Range("A2:A5").Copy($new.Range("A2:A5"))

# Header labels
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Fund code / name / size / position columns are text in the source data
# (e.g. fund codes keep leading zeros), so force text formatting before
# assigning the values.
$new.Range("B2:G5").NumberFormat = "@"

# Row 2: 011138 - 广发聚鸿六个月持有期混合A
$new.Range("B2").Value = "011138"
$new.Range("C2").Value = "广发聚鸿六个月持有期混合A"
$new.Range("D2").Value = "12.07"
$new.Range("E2").Value = "93.49"
$new.Range("F2").Value = "2.79"
$new.Range("G2").Value = "0.3368"
$new.Range("H2").Value = 10

# Row 3: 011140 - 广发聚鸿六个月持有期混合E
$new.Range("B3").Value = "011140"
$new.Range("C3").Value = "广发聚鸿六个月持有期混合E"
$new.Range("D3").Value = "12.07"
$new.Range("E3").Value = "93.49"
$new.Range("F3").Value = "2.79"
$new.Range("G3").Value = "0.3368"
$new.Range("H3").Value = 10

# Row 4: 162720 - 广发创业板两年定期开放混合
$new.Range("B4").Value = "162720"
$new.Range("C4").Value = "广发创业板两年定期开放混合"
$new.Range("D4").Value = "8.96"
$new.Range("E4").Value = "93.33"
$new.Range("F4").Value = "3.21"
$new.Range("G4").Value = "0.2876"
$new.Range("H4").Value = 9

# Row 5: 011139 - 广发聚鸿六个月持有期混合C
$new.Range("B5").Value = "011139"
$new.Range("C5").Value = "广发聚鸿六个月持有期混合C"
$new.Range("D5").Value = "0.64"
$new.Range("E5").Value = "93.49"
$new.Range("F5").Value = "2.79"
$new.Range("G5").Value = "0.0179"
$new.Range("H5").Value = 10

# ---------------------------------------------------------------------
# 2) Re-add "总计" at the end of the workbook and repopulate it with the
#    snapshot taken earlier plus the new "2022-Q1" row at the top.
# ---------------------------------------------------------------------
$tot2 = $wb.Worksheets.Add($null, $new)
$tot2.Name = "总计"

# Rebuild header + index-column styling from the "2021-Q4" sheet (same
# bold/centered header style, same bordered index-column style).
$q4.Range("B1:D1").Copy($tot2.Range("B1:D1"))
$q4.Range("A2:A7").Copy($tot2.Range("A2:A7"))

$tot2.Range("B1").Value = "日期"
$tot2.Range("C1").Value = "持有数量(只)"
$tot2.Range("D1").Value = "持有市值(亿元)"

$allRows = @(,@("2022-Q1", 4, 0.98)) + $existingTotals
for ($i = 0; $i -lt $allRows.Count; $i++) {
    $r = $i + 2
    $tot2.Cells.Item($r, 1).Value = $i
    $tot2.Cells.Item($r, 2).Value = $allRows[$i][0]
    $tot2.Cells.Item($r, 3).Value = $allRows[$i][1]
    $tot2.Cells.Item($r, 4).Value = $allRows[$i][2]
}
